# Appends, after the existing table (and the bookmark that wraps it), a
# horizontal-rule paragraph followed by an italic "Last updated on 5 Oct
# 2020" paragraph - mirroring the commit's addition of two new <w:p>
# elements just before the closing <w:sectPr/>.

$d = $word.ActiveDocument

# --- Step 1: get a real body paragraph after the table ---------------
# The document currently ends right at the table (which is wrapped by the
# bookmark pair). A collapsed Range at Content.End is ambiguous - it can
# resolve to "end of the last table cell" - so ask Word to insert a
# paragraph mark after that point first; that reliably lands outside the
# table, immediately after </w:bookmarkEnd>, which is where the new
# content belongs.
$endRange = $d.Range($d.Content.End, $d.Content.End)
$endRange.InsertParagraphAfter()

$newPara = $d.Paragraphs.Last
$insertPoint = $newPara.Range
$insertPoint.Collapse(1)

# --- Step 2: insert the horizontal rule + "last updated" paragraphs ---
# There is no higher-level Word object for a horizontal-line/"Insert
# Horizontal Line" shape or for redundant direct-formatting runs, so the
# new markup is supplied as literal WordOpenXML and dropped in at the
# insertion point via Range.InsertXML - the same mechanism Word itself
# uses under the hood for InsertXML/paste operations.
$newContentXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:o="urn:schemas-microsoft-com:office:office"><w:body><w:p><w:r><w:pict><v:rect style="width:0;height:1.5pt" o:hralign="center" o:hrstd="t" o:hr="t"/></w:pict></w:r></w:p><w:p><w:pPr><w:pStyle w:val="FirstParagraph"/></w:pPr><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t xml:space="preserve">Last updated on 5</w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t xml:space="preserve"> Oct 2020</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertPoint.InsertXML($newContentXml)

# --- Step 3: drop the leftover empty paragraph -------------------------
# InsertXML split the placeholder paragraph from step 1 around the
# insertion point, leaving its now-empty tail paragraph dangling after
# the new "last updated" paragraph (right before <w:sectPr/>). Delete it
# so the new paragraph is immediately followed by the section break,
# exactly like the target diff.
$paraCount = $d.Paragraphs.Count
$secondLastEnd = $d.Paragraphs.Item($paraCount - 1).Range.End
$lastEnd = $d.Paragraphs.Item($paraCount).Range.End
$cleanupRange = $d.Range($secondLastEnd - 1, $lastEnd)
$cleanupRange.Delete()
